$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Record"
$ws.Range("B16").Value = "Balanço Geral"
$ws.Range("C16").Value = "Trânsito"
$ws.Range("D16").Value = "2025-04-01T11:58"
$ws.Range("E16").Value = "Neutro"
$ws.Range("F16").Value = "Carreta invade o calçadão, derruba poste e deixa o Centro sem energia. Repórter *ao vivo*. Previsão é o caminhão ser retirado até 6h da tarde. Local isolado. Equipe da Enel no local. Guarda Municipal fez isolamento da área central para evitar déficit. Motorista é de Caxias do Sul (SC) e estava indo levar mercadoria para São Pedro da Aldeia. GPS estava marcando que aqui era uma rua. Estava escuro no momento. Entrevista com motorista da carreta e com comerciantes que ficaram sem energia. "
